# Update "provincias_spain" worksheet: refresh Cordoba/Cadiz province rows
# with corrected data and bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp: 13:46 -> 14:16
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 14:16"

# Rows whose city stays the same but the active/recovered/deaths counts change
$ws.Range("C14").Value = 48
$ws.Range("D14").Value = 731
$ws.Range("E14").Value = 40

$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 553
$ws.Range("E18").Value = 25

$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 511
$ws.Range("E20").Value = 18

$ws.Range("C33").Value = 5
$ws.Range("D33").Value = 297
$ws.Range("E33").Value = 14

# Rows 34-43: Cordoba and Cadiz move up into the sorted-by-total order with
# fresh figures, pushing the rows that used to occupy 34-39 down by two and
# leaving the remainder of the block re-aligned to the new sort order.
$ws.Range("A34").Value = "Cordoba"
$ws.Range("B34").Value = 291
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 285
$ws.Range("E34").Value = 6

$ws.Range("A35").Value = "Cadiz"
$ws.Range("B35").Value = 278
$ws.Range("C35").Value = 4
$ws.Range("D35").Value = 270
$ws.Range("E35").Value = 4

$ws.Range("A36").Value = "Segovia"
$ws.Range("B36").Value = 271
$ws.Range("C36").Value = 32
$ws.Range("D36").Value = 212
$ws.Range("E36").Value = 27

$ws.Range("A37").Value = "Castello/Castellon"
$ws.Range("B37").Value = 269
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 257
$ws.Range("E37").Value = 11

$ws.Range("A38").Value = "Guadalajara"
$ws.Range("B38").Value = 263
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 257
$ws.Range("E38").Value = 4

$ws.Range("A39").Value = "Badajoz"
$ws.Range("B39").Value = 257
$ws.Range("C39").Value = 5
$ws.Range("D39").Value = 248
$ws.Range("E39").Value = 4

$ws.Range("A40").Value = "Mallorca"
$ws.Range("B40").Value = 210
$ws.Range("C40").Value = 18
$ws.Range("D40").Value = 194
$ws.Range("E40").Value = 12

$ws.Range("A41").Value = "Avila"
$ws.Range("B41").Value = 201
$ws.Range("C41").Value = 23
$ws.Range("D41").Value = 163
$ws.Range("E41").Value = 15

$ws.Range("A42").Value = "Ourense"
$ws.Range("B42").Value = 189
$ws.Range("C42").Value = 25
$ws.Range("D42").Value = 186
$ws.Range("E42").Value = 3

$ws.Range("A43").Value = "Soria"
$ws.Range("B43").Value = 179
$ws.Range("C43").Value = 14
$ws.Range("D43").Value = 152
$ws.Range("E43").Value = 13

# Two more same-city, corrected-metrics rows
$ws.Range("C46").Value = 5
$ws.Range("D46").Value = 105
$ws.Range("E46").Value = 5

$ws.Range("C49").Value = 2
$ws.Range("D49").Value = 74
$ws.Range("E49").Value = 1
